# [Fix]: exclusion of 8 redundant metrics
#
# Eight redundant metrics (MBRAE, UMBRAE, STDAPE, RMSPE, MRE, MRAE, MDRAE,
# GMRAE) are removed from both worksheets ("LMN" and "STS"). Column A (the
# numeric row id 0..33) is left untouched; columns B (metric name) and C
# (metric value) for the surviving metrics are compacted upward to fill
# rows 2..26, and the now-unused trailing rows (27..34) are deleted.

$wb = $excel.ActiveWorkbook

# Ordered (name, sheet1 value, sheet2 value) for every metric row, in the
# original top-to-bottom order (this is the full pre-edit metric list).
$metrics = @(
    @{ Name = "MAPE";   S1 =  0.2222779178827255;  S2 =  0.1242740884687008 },
    @{ Name = "WAPE";   S1 =  0.06382605258508071; S2 =  0.0991655832188242 },
    @{ Name = "MAE";    S1 =  0.01616899426102431; S2 =  0.07637765802491184 },
    @{ Name = "MAAPE";  S1 =  0.1675605089206508;  S2 =  0.1211394304000615 },
    @{ Name = "MASE";   S1 =  0.05024249454497252; S2 =  0.08235761296778825 },
    @{ Name = "MSE";    S1 =  0.001609146109663624;S2 =  0.02448073665962463 },
    @{ Name = "RMSE";   S1 =  0.04011416345461567; S2 =  0.1564632118410735 },
    @{ Name = "NRMSE";  S1 =  0.04614367054262412; S2 =  0.04898798395730419 },
    @{ Name = "R^2";    S1 =  0.982471647529276;   S2 =  0.9670516720568563 },
    @{ Name = "Pearson";S1 =  0.9912020044369025;  S2 =  0.9837791698512718 },
    @{ Name = "MBRAE";  S1 =  0.1290307591529484;  S2 =  0.1408164153293056; Remove = $true },
    @{ Name = "UMBRAE"; S1 =  0.1481461722201131;  S2 =  0.1638956072272696; Remove = $true },
    @{ Name = "ME";     S1 = -0.0008396788340461281;S2 = 0.01741963747947288 },
    @{ Name = "MAD";    S1 =  0.01616899426102431; S2 =  0.07637765802491184 },
    @{ Name = "GMAE";   S1 =  0.006223043926114697;S2 =  0.02966695604147439 },
    @{ Name = "MDAE";   S1 =  0.00640214445590318; S2 =  0.02784852562745999 },
    @{ Name = "MPE";    S1 = -0.02591341006578651; S2 = -0.01640369233933391 },
    @{ Name = "MDAPE";  S1 =  0.0855192477886033;  S2 =  0.08396310123967905 },
    @{ Name = "SMAPE";  S1 =  0.2298921048055551;  S2 =  0.1236855444827733 },
    @{ Name = "SMDAPE"; S1 =  0.08610884874931543; S2 =  0.08576213714464509 },
    @{ Name = "STDAE";  S1 =  0.04359966572510861; S2 =  0.1664328558987114 },
    @{ Name = "STDAPE"; S1 =  0.6219497845759493;  S2 =  0.2229077734066728; Remove = $true },
    @{ Name = "RMSPE";  S1 =  0.5703016604046965;  S2 =  0.1734451196201527; Remove = $true },
    @{ Name = "RMDSPE"; S1 =  0.0855192477886033;  S2 =  0.08396310123967905 },
    @{ Name = "RMSSE";  S1 =  0.1246481757621074;  S2 =  0.1687134297872976 },
    @{ Name = "INRSE";  S1 =  0.1323946844504115;  S2 =  0.1815167428727822 },
    @{ Name = "RRSE";   S1 =  0.1323946844504115;  S2 =  0.1815167428727822 },
    @{ Name = "MRE";    S1 = -0.8834949954528452;  S2 = -0.3704618702239451; Remove = $true },
    @{ Name = "RAE";    S1 =  0.05890744415776347; S2 =  0.1110060452924479 },
    @{ Name = "MRAE";   S1 =  1.102620171988271;   S2 =  0.6840012904260514; Remove = $true },
    @{ Name = "MDRAE";  S1 =  0.04893006398268832; S2 =  0.08375497011968289; Remove = $true },
    @{ Name = "GMRAE";  S1 =  0.05037199315190379; S2 =  0.06912027758254025; Remove = $true },
    @{ Name = "MDA";    S1 =  0.9494949494949495;  S2 =  0.9494949494949495 }
)

$kept = $metrics | Where-Object { -not $_.Remove }

foreach ($ws in $wb.Worksheets) {
    $valueField = if ($ws.Name -eq "LMN") { "S1" } else { "S2" }

    $row = 2
    foreach ($m in $kept) {
        $ws.Cells.Item($row, 2).Value = $m.Name
        $ws.Cells.Item($row, 3).Value = $m[$valueField]
        $row = $row + 1
    }

    # $row is now one past the last surviving data row (row 27 if 25
    # metrics survive, i.e. rows 2..26 filled) -> delete the remaining
    # now-stale rows down through the old last row (34), bottom-up.
    for ($r = 34; $r -ge $row; $r--) {
        $ws.Rows.Item($r).Delete()
    }
}
